$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 1.700000047683716, 2.720000028610229, 3.109999895095825, 1.700000047683716, 1220418767),
    @(3, 3.403000116348267, 3.739000082015991, 4.025000095367432, 2.849999904632568, 1220418767),
    @(4, 3.536999940872192, 3.221999883651733, 3.928999900817871, 2.971999883651733, 1220418767),
    @(5, 2.54200005531311, 2.322000026702881, 2.650000095367432, 1.848000049591064, 1220418767),
    @(6, 2.788000106811523, 3.184999942779541, 3.239000082015991, 2.756500005722046, 1220418767),
    @(7, 3.075000047683716, 3.427000045776367, 3.476000070571899, 3, 1220418767),
    @(8, 4.323999881744385, 4.144999980926514, 4.5, 4.014999866485596, 1220418767),
    @(9, 4.328000068664551, 5.081999778747559, 5.199999809265137, 4.213500022888184, 1220418767),
    @(10, 6.815999984741211, 7.59499979019165, 7.818999767303467, 6.722000122070312, 1220418767),
    @(11, 8.715000152587891, 9.237000465393066, 9.680000305175779, 8.479999542236328, 1220418767),
    @(12, 11.69999980926514, 9.949000358581545, 12.06900024414062, 8.935000419616699, 1220418767),
    @(13, 10.15900039672852, 12.79199981689453, 13.04699993133545, 10.10200023651123, 1220418767),
    @(14, 12.4709997177124, 13.36299991607666, 13.44299983978272, 11.25, 1220418767),
    @(15, 14.27099990844727, 13.82100009918213, 17.6599006652832, 13.32610034942627, 1220418767),
    @(16, 16.64399909973145, 13.8149995803833, 16.68600082397461, 12.20499992370606, 1220418767),
    @(17, 13.39999961853027, 16.84700012207031, 16.95999908447266, 12.94799995422363, 1220418767),
    @(18, 20.8430004119873, 24.35300064086914, 24.75, 19.03800010681152, 1220418767),
    @(19, 30.86000061035156, 31.78800010681152, 34.08359909057617, 30.03300094604492, 1220418767),
    @(20, 31.16600036621094, 31.35700035095215, 34.94400024414062, 29.10600090026855, 1220418767),
    @(21, 40.39899826049805, 46.56600189208984, 48.2869987487793, 39.54199981689453, 1220418767),
    @(22, 40.29899978637695, 63.22900009155273, 66.57399749755859, 33.45500183105469, 1220418767),
    @(23, 94.80100250244141, 102.4000015258789, 110.7919998168945, 88.31800079345703, 1220418767),
    @(24, 103.0630035400391, 92.54299926757812, 113, 90.62200164794922, 1220418767),
    @(25, 113.7990036010742, 109.859001159668, 124.745002746582, 107.6050033569336, 1220418767),
    @(26, 115, 118.2509994506836, 130.1609954833984, 108.7259979248047, 1220418767),
    @(27, 146.0980072021484, 149.9929962158203, 165, 139.3009948730469, 1220418767),
    @(28, 136.2460021972656, 146.6730041503906, 150.6719970703125, 128.5899963378906, 1220418767),
    @(29, 139.1900024414062, 96.42400360107422, 139.3670043945312, 78, 1220418767),
    @(30, 68.75599670410156, 42.68199920654297, 74.13400268554688, 41.11700057983398, 1220418767),
    @(31, 32.08000183105469, 34.83000183105469, 40.68000030517578, 29.71999931335449, 1220418767),
    @(32, 27.29999923706055, 34.22999954223633, 36.15999984741211, 23.6299991607666, 1220418767),
    @(33, 35.68999862670898, 49.27000045776367, 50.16999816894531, 34.65000152587891, 1220418767),
    @(34, 47.25, 48.45000076293945, 49.86000061035156, 43.5, 1220418767),
    @(35, 65.04000091552734, 67.58000183105469, 71.43000030517578, 60.5099983215332, 1220418767),
    @(36, 54.22000122070312, 47.18999862670898, 56.19599914550781, 45.5, 1220418767),
    @(37, 76.44000244140625, 80.06999969482422, 83.57499694824219, 70.61499786376953, 1220418767),
    @(38, 77.34999847412109, 70.19999694824219, 79.29000091552734, 68.01499938964844, 1220418767),
    @(39, 65.95999908447266, 61.20000076293945, 70.04000091552734, 58.04000091552734, 1220418767),
    @(40, 80, 78.20999908447266, 84.37000274658203, 77.12000274658203, 1220418767),
    @(41, 107.5, 116.8000030517578, 122, 99.0500030517578, 1220418767),
    @(42, 95.8000030517578, 95, 101.4499969482422, 69.83999633789062, 1220418767),
    @(43, 114.3099975585938, 122.2099990844727, 129.2899932861328, 111.3300018310547, 1220418767),
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 4).Value2 = $item[1]
    $ws.Cells.Item($r, 5).Value2 = $item[2]
    $ws.Cells.Item($r, 6).Value2 = $item[3]
    $ws.Cells.Item($r, 7).Value2 = $item[4]
    $ws.Cells.Item($r, 8).Value2 = $item[5]
    $ws.Cells.Item($r, 9).Value = "SHOP"
}
